$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three existing August rows with revised totals
$ws.Range("B2").Value = 23923.99
$ws.Range("B3").Value = 30379.91
$ws.Range("B4").Value = 15076.49

# Insert a new row for August (Dia 6), pushing all following rows down by one
$ws.Rows("5:5").Insert()

# Populate the newly inserted row
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 16426.77
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 2025
$ws.Range("E5").Value = "08/2025"
